$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell E8 text from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 on the active sheet to match the sheetView selection in the diff
$ws.Activate()
$ws.Range("E8").Select()
